$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Value = 5500.8832938526803
